$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header relabeling (D1:BQ1) ---
$ws.Range('D1').Value2 = 'Index child (f) log(HR)'
$ws.Range('E1').Value2 = 'SE'
$ws.Range('F1').Value2 = 'n'
$ws.Range('G1').Value2 = 'Index child (m) log(HR)'
$ws.Range('H1').Value2 = 'SE'
$ws.Range('I1').Value2 = 'n'
$ws.Range('J1').Value2 = 'Brother log(HR)'
$ws.Range('K1').Value2 = 'SE'
$ws.Range('L1').Value2 = 'n'
$ws.Range('M1').Value2 = 'Sister log(HR)'
$ws.Range('N1').Value2 = 'SE'
$ws.Range('O1').Value2 = 'n'
$ws.Range('P1').Value2 = 'Mat. half sister log(HR)'
$ws.Range('Q1').Value2 = 'SE'
$ws.Range('R1').Value2 = 'n'
$ws.Range('S1').Value2 = 'Pat. half sister log(HR)'
$ws.Range('T1').Value2 = 'SE'
$ws.Range('U1').Value2 = 'n'
$ws.Range('V1').Value2 = 'Mat. half brother log(HR)'
$ws.Range('W1').Value2 = 'SE'
$ws.Range('X1').Value2 = 'n'
$ws.Range('Y1').Value2 = 'Pat. half brother log(HR)'
$ws.Range('Z1').Value2 = 'SE'
$ws.Range('AA1').Value2 = 'n'
$ws.Range('AB1').Value2 = 'Mother log(HR)'
$ws.Range('AC1').Value2 = 'SE'
$ws.Range('AD1').Value2 = 'n'
$ws.Range('AE1').Value2 = 'Father log(HR)'
$ws.Range('AF1').Value2 = 'SE'
$ws.Range('AG1').Value2 = 'n'
$ws.Range('AH1').Value2 = 'Mat. grandmother log(HR)'
$ws.Range('AI1').Value2 = 'SE'
$ws.Range('AJ1').Value2 = 'n'
$ws.Range('AK1').Value2 = 'Mat. grandfather log(HR)'
$ws.Range('AL1').Value2 = 'SE'
$ws.Range('AM1').Value2 = 'n'
$ws.Range('AN1').Value2 = 'Pat. grandmother log(HR)'
$ws.Range('AO1').Value2 = 'SE'
$ws.Range('AP1').Value2 = 'n'
$ws.Range('AQ1').Value2 = 'Pat. grandfather log(HR)'
$ws.Range('AR1').Value2 = 'SE'
$ws.Range('AS1').Value2 = 'n'
$ws.Range('AT1').Value2 = 'Mat. aunt log(HR)'
$ws.Range('AU1').Value2 = 'SE'
$ws.Range('AV1').Value2 = 'n'
$ws.Range('AW1').Value2 = 'Mat. uncle log(HR)'
$ws.Range('AX1').Value2 = 'SE'
$ws.Range('AY1').Value2 = 'n'
$ws.Range('AZ1').Value2 = 'Pat. aunt log(HR)'
$ws.Range('BA1').Value2 = 'SE'
$ws.Range('BB1').Value2 = 'n'
$ws.Range('BC1').Value2 = 'Pat. uncle log(HR)'
$ws.Range('BD1').Value2 = 'SE'
$ws.Range('BE1').Value2 = 'n'
$ws.Range('BF1').Value2 = 'Mat. cousin (f) log(HR)'
$ws.Range('BG1').Value2 = 'SE'
$ws.Range('BH1').Value2 = 'n'
$ws.Range('BI1').Value2 = 'Mat. cousin (m) log(HR)'
$ws.Range('BJ1').Value2 = 'SE'
$ws.Range('BK1').Value2 = 'n'
$ws.Range('BL1').Value2 = 'Pat. cousin (f) log(HR)'
$ws.Range('BM1').Value2 = 'SE'
$ws.Range('BN1').Value2 = 'n'
$ws.Range('BO1').Value2 = 'Pat. cousin (m) log(HR)'
$ws.Range('BP1').Value2 = 'SE'
$ws.Range('BQ1').Value2 = 'n'

# --- Column C diagnosis relabeling (C2:C93) ---
$ws.Range('C2').Value2 = 'ASD'
$ws.Range('C3').Value2 = 'Any mental'
$ws.Range('C4').Value2 = 'Organic mental'
$ws.Range('C5').Value2 = 'Psychoactive sub use'
$ws.Range('C6').Value2 = 'Schizophrenia spectrum '
$ws.Range('C7').Value2 = 'Schizophrenia'
$ws.Range('C8').Value2 = 'Any mood '
$ws.Range('C9').Value2 = 'Bipolar disorder'
$ws.Range('C10').Value2 = 'Depression'
$ws.Range('C11').Value2 = 'Neurotic/stress disorder'
$ws.Range('C12').Value2 = 'OCD'
$ws.Range('C13').Value2 = 'Behav synd-physiol'
$ws.Range('C14').Value2 = 'Anorexia nervosa'
$ws.Range('C15').Value2 = 'Adult personality disorder'
$ws.Range('C16').Value2 = 'Intellectual disability'
$ws.Range('C17').Value2 = 'Psych dev dis.not ASD'
$ws.Range('C18').Value2 = 'Behav dis-child onset'
$ws.Range('C19').Value2 = 'ADHD'
$ws.Range('C20').Value2 = 'Tic disorder'
$ws.Range('C21').Value2 = 'Mental-unspecified'
$ws.Range('C22').Value2 = 'Any diabetes'
$ws.Range('C23').Value2 = 'Diabetes outside preg'
$ws.Range('C24').Value2 = 'Type 1 diabetes'
$ws.Range('C25').Value2 = 'Type 2 diabetes'
$ws.Range('C26').Value2 = 'Chronic+gest diab'
$ws.Range('C27').Value2 = 'Obesity'
$ws.Range('C28').Value2 = 'Any hypertension'
$ws.Range('C29').Value2 = 'Hyper outside preg'
$ws.Range('C30').Value2 = 'Chronic+gest hyper'
$ws.Range('C31').Value2 = 'Preeclam/eclam'
$ws.Range('C32').Value2 = 'Any birth defect'
$ws.Range('C33').Value2 = 'CNS'
$ws.Range('C34').Value2 = 'Eye'
$ws.Range('C35').Value2 = 'Ear'
$ws.Range('C36').Value2 = 'Heart'
$ws.Range('C37').Value2 = 'Respiratory'
$ws.Range('C38').Value2 = 'Lip'
$ws.Range('C39').Value2 = 'Digestive system'
$ws.Range('C40').Value2 = 'Genital'
$ws.Range('C41').Value2 = 'Urinary tract'
$ws.Range('C42').Value2 = 'Musculoskeletal'
$ws.Range('C43').Value2 = 'Skin'
$ws.Range('C44').Value2 = 'Other/chromos'
$ws.Range('C45').Value2 = 'Chro/gene dis_ASD spe'
$ws.Range('C46').Value2 = 'Any neurologic'
$ws.Range('C47').Value2 = 'Inflammatory of CNS'
$ws.Range('C48').Value2 = 'Systemic atrophies'
$ws.Range('C49').Value2 = 'Extrapyramid'
$ws.Range('C50').Value2 = 'Other degenerative '
$ws.Range('C51').Value2 = 'Demyelinating of CNS'
$ws.Range('C52').Value2 = 'Episodic'
$ws.Range('C53').Value2 = 'Epilepsy'
$ws.Range('C54').Value2 = 'Nerve disorder'
$ws.Range('C55').Value2 = 'Polyneuropath'
$ws.Range('C56').Value2 = 'Myoneural '
$ws.Range('C57').Value2 = 'Cerebral palsy '
$ws.Range('C58').Value2 = 'Other neurologic'
$ws.Range('C59').Value2 = 'Type 1 diabetes'
$ws.Range('C60').Value2 = 'Thyrotoxicosis'
$ws.Range('C61').Value2 = 'Thyroiditis'
$ws.Range('C62').Value2 = 'Pri adrenocortical '
$ws.Range('C63').Value2 = 'Rheumatoid arthritis'
$ws.Range('C64').Value2 = 'Juvenile arthritis'
$ws.Range('C65').Value2 = 'Dermatopolymyositis'
$ws.Range('C66').Value2 = 'Polymyalgia'
$ws.Range('C67').Value2 = 'Scleroderma'
$ws.Range('C68').Value2 = 'Lupus erythema'
$ws.Range('C69').Value2 = 'Sjogren'
$ws.Range('C70').Value2 = 'Ankylos spondil '
$ws.Range('C71').Value2 = 'Granulomato'
$ws.Range('C72').Value2 = 'Celiac'
$ws.Range('C73').Value2 = 'Crohn'
$ws.Range('C74').Value2 = 'Ulcerative colitis'
$ws.Range('C75').Value2 = 'Pernicious anem'
$ws.Range('C76').Value2 = 'Hemolytic anem'
$ws.Range('C77').Value2 = 'Purpura'
$ws.Range('C78').Value2 = 'Multiple sclerosis'
$ws.Range('C79').Value2 = 'Guillain-Bar'
$ws.Range('C80').Value2 = 'Myasthen grav'
$ws.Range('C81').Value2 = 'Pemphigus'
$ws.Range('C82').Value2 = 'Psoriasis'
$ws.Range('C83').Value2 = 'Alopecia areata'
$ws.Range('C84').Value2 = 'Vitiligo'
$ws.Range('C85').Value2 = 'Any endocrine '
$ws.Range('C86').Value2 = 'Any connective'
$ws.Range('C87').Value2 = 'Any gastrointest'
$ws.Range('C88').Value2 = 'Any blood'
$ws.Range('C89').Value2 = 'Any nervous'
$ws.Range('C90').Value2 = 'Any skin'
$ws.Range('C91').Value2 = 'Any autoimmune'
$ws.Range('C92').Value2 = 'Asthma'
$ws.Range('C93').Value2 = 'Allergies'
